$d = $word.ActiveDocument

# Remove the existing "_GoBack" bookmark (it currently sits after the
# "資訊機器人沒電時..." paragraph and needs to move).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Find the "聊時事" substring inside the "這台機器人能跟他聊時事、..." sentence
# and wrap it with a new "_GoBack" bookmark, which splits the run into three
# pieces: "這台機器人能跟他" | "聊時事" | "、幫助他解決生活中的疑難雜症".
$rng = $d.Content
$found = $rng.Find.Execute("聊時事", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $d.Bookmarks.Add("_GoBack", $rng)
}
